$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H96").Value = 715.6667
$ws.Range("I96").Value = 299.33334
$ws.Range("J96").Value = 1132
$ws.Range("K96").Value = 898.0000200000001
$ws.Range("L96").Value = 3396
$ws.Range("M96").Value = 474.9999799999999
$ws.Range("N96").ClearContents()

$ws.Range("H112").Value = 1819.8
$ws.Range("I112").Value = 1700
$ws.Range("J112").Value = 1999.5
$ws.Range("K112").Value = 5100
$ws.Range("L112").Value = 5998.5
$ws.Range("M112").Value = -3992
$ws.Range("N112").Value = -8214.5

$ws.Range("H137").Value = 2841.1428
$ws.Range("I137").Value = 2841.1428
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 8523.428400000001
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -5973.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 962
$ws.Range("I2").Value = 962
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 962
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -849

$ws.Range("H22").Value = 3000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 3000
$ws.Range("N22").Value = -3598

$ws.Range("H61").Value = 10475.154
$ws.Range("I61").Value = 12954.889
$ws.Range("J61").Value = 4895.75
$ws.Range("K61").Value = 12954.889
$ws.Range("L61").Value = 4895.75
$ws.Range("M61").Value = -12742.889
$ws.Range("N61").Value = -5319.75

$ws.Range("H112").Value = 27000
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 27000
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 27000
$ws.Range("N112").Value = -29954

$ws.Range("H116").Value = 962
$ws.Range("I116").Value = 962
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 962
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1332

$ws.Range("H136").Value = 10475.154
$ws.Range("I136").Value = 12954.889
$ws.Range("J136").Value = 4895.75
$ws.Range("K136").Value = 38864.667
$ws.Range("L136").Value = 14687.25
$ws.Range("M136").Value = -36314.667
$ws.Range("N136").Value = -19787.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 962
$ws.Range("I3").Value = 962
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 962
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -848

$ws.Range("H81").Value = 35111.6
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 35111.6
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 35111.6
$ws.Range("N81").Value = -37233.6

$ws.Range("H82").Value = 23697.143
$ws.Range("I82").Value = 5980
$ws.Range("J82").Value = 130000
$ws.Range("K82").Value = 5980
$ws.Range("L82").Value = 130000
$ws.Range("M82").Value = -5597
$ws.Range("N82").ClearContents()

$ws.Range("H84").Value = 35111.6
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 35111.6
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 105334.8
$ws.Range("N84").Value = -115942.8

$ws.Range("H85").Value = 23697.143
$ws.Range("I85").Value = 5980
$ws.Range("J85").Value = 130000
$ws.Range("K85").Value = 5980
$ws.Range("L85").Value = 130000
$ws.Range("M85").Value = -4654
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 30.333334
$ws.Range("I7").Value = 10.166667
$ws.Range("J7").Value = 70.666664
$ws.Range("K7").Value = 10.166667
$ws.Range("L7").Value = 70.666664
$ws.Range("M7").Value = 102.833333
$ws.Range("N7").Value = -296.666664

$ws.Range("H16").Value = 4183.3335
$ws.Range("I16").Value = 1650
$ws.Range("J16").Value = 5450
$ws.Range("K16").Value = 1650
$ws.Range("L16").Value = 5450
$ws.Range("M16").Value = -1363
$ws.Range("N16").Value = -6024

$ws.Range("H22").Value = 805.44446
$ws.Range("I22").Value = 793.625
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 793.625
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -443.625
$ws.Range("N22").ClearContents()

$ws.Range("H31").Value = 3995.2
$ws.Range("I31").Value = 2817.8
$ws.Range("J31").Value = 5172.6
$ws.Range("K31").Value = 2817.8
$ws.Range("L31").Value = 5172.6
$ws.Range("M31").Value = -2522.8
$ws.Range("N31").ClearContents()

$ws.Range("H34").Value = 3995.2
$ws.Range("I34").Value = 2817.8
$ws.Range("J34").Value = 5172.6
$ws.Range("K34").Value = 2817.8
$ws.Range("L34").Value = 5172.6
$ws.Range("M34").Value = -2615.8
$ws.Range("N34").ClearContents()

$ws.Range("H103").Value = 14532.667
$ws.Range("I103").Value = 14532.667
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 14532.667
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -13360.667

$ws.Range("H113").Value = 4183.3335
$ws.Range("I113").Value = 1650
$ws.Range("J113").Value = 5450
$ws.Range("K113").Value = 1650
$ws.Range("L113").Value = 5450
$ws.Range("M113").Value = 520
$ws.Range("N113").Value = -9790

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 5105
$ws.Range("I38").Value = 5506.6665
$ws.Range("J38").Value = 4502.5
$ws.Range("K38").Value = 16519.9995
$ws.Range("L38").Value = 13507.5
$ws.Range("M38").Value = -16172.9995
$ws.Range("N38").Value = -14201.5

$ws.Range("H132").Value = 998
$ws.Range("I132").Value = 995
$ws.Range("J132").Value = 998.75
$ws.Range("K132").Value = 8955
$ws.Range("L132").Value = 8988.75
$ws.Range("M132").Value = -6425
$ws.Range("N132").Value = -14048.75

$ws.Range("H134").Value = 4900
$ws.Range("I134").Value = 4900
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 14700
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -9630

$ws.Range("H137").Value = 3000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 9000
$ws.Range("N137").Value = -19200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()

$ws.Range("H102").Value = 699.5
$ws.Range("I102").Value = 599
$ws.Range("J102").Value = 800
$ws.Range("K102").Value = 599
$ws.Range("L102").Value = 800
$ws.Range("M102").Value = 1023
$ws.Range("N102").ClearContents()

$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6550

$ws.Range("H132").Value = 3214.6
$ws.Range("I132").Value = 2641
$ws.Range("J132").Value = 4075
$ws.Range("K132").Value = 7923
$ws.Range("L132").Value = 12225
$ws.Range("M132").Value = -5393
$ws.Range("N132").Value = -17285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 5430.8887
$ws.Range("I132").Value = 4976.2
$ws.Range("J132").Value = 5999.25
$ws.Range("K132").Value = 14928.6
$ws.Range("L132").Value = 17997.75
$ws.Range("M132").Value = -12398.6
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 1175.5
$ws.Range("I6").Value = 700
$ws.Range("J6").Value = 2602
$ws.Range("K6").Value = 700
$ws.Range("L6").Value = 2602
$ws.Range("M6").Value = -585
$ws.Range("N6").Value = -2832
